$d = $word.ActiveDocument

# 1) Capitalize "g" -> "G" and, in the same operation, split the single
#    paragraph into five paragraphs by inserting paragraph marks (^p).
#    This pushes the trailing bookmark (originally sharing the "g"
#    paragraph) down into its own, final, empty paragraph - exactly like
#    the diff shows.
$d.Content.Find.Execute("g", $true, $false, $false, $false, $false, $true, 1, $false, `
    "G^pBagles ^pHot sauce^pOranges^p", 2)

# 2) The "Bagles " paragraph needs proofErr spell-check markers around
#    the word "Bagles". Replace that paragraph's content via InsertXML so
#    we can place <w:proofErr> elements precisely.
$p2 = $d.Paragraphs(2)
$r = $p2.Range

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Bagles</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
